$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings are kept as literal text (matches source data format)
$textCells = @('D4', 'D5', 'D6', 'D7', 'D9', 'D10', 'D13', 'D14', 'D15', 'D16', 'D17', 'D18', 'D19', 'D20', 'D23', 'D24', 'D25', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

# Apply updated cell values from the crypto price refresh
$ws.Range('D2').Value = '27.083.11'
$ws.Range('E2').Value = '  +0.12%  '
$ws.Range('D3').Value = '1.818.32'
$ws.Range('E3').Value = '  -0.43%  '
$ws.Range('D4').Value = '0.9990'
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').Value = '311.63'
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('D6').Value = '0.9988'
$ws.Range('E6').Value = '  -0.21%  '
$ws.Range('D7').Value = '0.4446'
$ws.Range('E7').Value = '  +5.06%  '
$ws.Range('E8').Value = '  +1.82%  '
$ws.Range('D9').Value = '0.07436'
$ws.Range('E9').Value = '  +2.93%  '
$ws.Range('D10').Value = '0.8716'
$ws.Range('E10').Value = '  +3.29%  '
$ws.Range('E11').Value = '  +0.45%  '
$ws.Range('D12').Value = '1.818.11'
$ws.Range('E12').Value = '  -0.41%  '
$ws.Range('D13').Value = '6.710'
$ws.Range('E13').Value = '  +0.85%  '
$ws.Range('D14').Value = '93.86'
$ws.Range('E14').Value = '  +4.54%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = '5.328'
$ws.Range('E15').Value = '  +0.69%  '
$ws.Range('B16').Value = 'TRON'
$ws.Range('C16').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D16').Value = '0.07100'
$ws.Range('E16').Value = '  +0.79%  '
$ws.Range('D17').Value = '0.9995'
$ws.Range('E17').Value = '  -0.32%  '
$ws.Range('D18').Value = '0.000008748'
$ws.Range('E18').Value = '  +0.07%  '
$ws.Range('D19').Value = '0.9992'
$ws.Range('E19').Value = '  -0.13%  '
$ws.Range('D20').Value = '14.98'
$ws.Range('E20').Value = '  +0.70%  '
$ws.Range('D21').Value = '27.104.96'
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('E22').Value = '  +1.46%  '
$ws.Range('D23').Value = '10.93'
$ws.Range('E23').Value = '  +0.99%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').Value = '1.979'
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '151.42'
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D26').Value = '2.250'
$ws.Range('E26').Value = '  -0.16%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '18.49'
$ws.Range('E27').Value = '  +1.50%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = '5.293'
$ws.Range('E28').Value = '  +0.83%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').Value = '118.12'
$ws.Range('E29').Value = '  +1.28%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').Value = '0.08830'
$ws.Range('E30').Value = '  +1.49%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = '0.7631'
$ws.Range('E31').Value = '  +3.42%  '
$ws.Range('B32').Value = 'ARBITRUM'
$ws.Range('C32').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D32').Value = '1.172'
$ws.Range('E32').Value = '  -0.52%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '4.559'
$ws.Range('E33').Value = '  +3.11%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = '2.884'
$ws.Range('E34').Value = '  -0.47%  '
$ws.Range('B35').Value = 'Frax'
$ws.Range('C35').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D35').Value = '0.9985'
$ws.Range('E35').Value = '  -0.21%  '
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').Value = '1.101'
$ws.Range('E36').Value = '  +0.76%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '0.01981'
$ws.Range('E37').Value = '  +1.90%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = '0.05260'
$ws.Range('E38').Value = '  +0.25%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').Value = '7.278'
$ws.Range('E39').Value = '  -0.56%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '0.5289'
$ws.Range('E40').Value = '  +4.42%  '
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').Value = '2.856'
$ws.Range('E41').Value = '  -0.57%  '
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').Value = '0.1711'
$ws.Range('E42').Value = '  +1.53%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = '2.152'
$ws.Range('E43').Value = '  +9.41%  '
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').Value = '8.696'
$ws.Range('E44').Value = '  +1.69%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = '0.5028'
$ws.Range('E45').Value = '  +6.40%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '10.59'
$ws.Range('E46').Value = '  +0.72%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').Value = '1.705'
$ws.Range('E47').Value = '  +3.23%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = '105.08'
$ws.Range('E48').Value = '  -0.76%  '
$ws.Range('B49').Value = 'PaxDollar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D49').Value = '0.9982'
$ws.Range('E49').Value = '  -0.21%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.06360'
$ws.Range('E50').Value = '  +0.53%  '
$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D51').Value = '0.9303'
$ws.Range('E51').Value = '  +3.02%  '
